# Insert a new data row at row 506 of Sheet1 (Espinaca - Femacal de La Calera).
# This shifts the existing rows 506..571 down to 507..572 and fills the
# newly freed row 506 with a new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 506 downward (this also extends the used range / dimension).
$ws.Rows.Item(506).Insert()

# Populate the new row 506 with the new record's values.
$ws.Range("A506").Value = 3
$ws.Range("B506").Value = "Femacal de La Calera"
$ws.Range("C506").Value = "Coquimbo"
$ws.Range("D506").Value = 45124
$ws.Range("E506").Value = 5
$ws.Range("F506").Value = 100112012
$ws.Range("G506").Value = "Espinaca"
$ws.Range("H506").Value = "Sin especificar"
$ws.Range("I506").Value = "Primera"
$ws.Range("J506").Value = 130
$ws.Range("K506").Value = 4000
$ws.Range("L506").Value = 4800
$ws.Range("M506").Value = 4400
$ws.Range("N506").Value = "$/docena de atados (3 kilos)"
$ws.Range("O506").Value = "Provincia de Quillota"
$ws.Range("P506").Value = 1467
$ws.Range("Q506").Value = 3
$ws.Range("R506").Value = "Hortaliza"
